$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data cells while preserving their original
# (default / unstyled) cell style and forcing the values to remain plain
# text, since the source data uses text-formatted numbers (e.g. "29.471.95",
# "1.000") that must not be reinterpreted as numeric values by Excel.
function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '29.471.95'
Set-TextValue $ws.Range('E2') '  +1.74%  '
Set-TextValue $ws.Range('D3') '1.858.21'
Set-TextValue $ws.Range('E3') '  +1.25%  '
Set-TextValue $ws.Range('D4') '0.9994'
Set-TextValue $ws.Range('E4') '  -0.01%  '
Set-TextValue $ws.Range('D5') '245.37'
Set-TextValue $ws.Range('E5') '  -0.09%  '
Set-TextValue $ws.Range('D6') '0.6947'
Set-TextValue $ws.Range('E6') '  +0.68%  '
Set-TextValue $ws.Range('D7') '1.000'
Set-TextValue $ws.Range('E7') '  +0.02%  '
Set-TextValue $ws.Range('B8') 'Dogecoin'
Set-TextValue $ws.Range('C8') 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D8') '0.07675'
Set-TextValue $ws.Range('E8') '  -0.54%  '
Set-TextValue $ws.Range('B9') 'Cardano'
Set-TextValue $ws.Range('C9') 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range('D9') '0.3064'
Set-TextValue $ws.Range('E9') '  +0.29%  '
Set-TextValue $ws.Range('D10') '23.61'
Set-TextValue $ws.Range('E10') '  +0.24%  '
Set-TextValue $ws.Range('E11') '  -0.44%  '
Set-TextValue $ws.Range('D12') '5.149'
Set-TextValue $ws.Range('E12') '  +1.32%  '
Set-TextValue $ws.Range('D13') '1.856.39'
Set-TextValue $ws.Range('E13') '  +1.13%  '
Set-TextValue $ws.Range('D14') '91.15'
Set-TextValue $ws.Range('E14') '  +0.71%  '
Set-TextValue $ws.Range('D15') '0.6917'
Set-TextValue $ws.Range('E15') '  +1.57%  '
Set-TextValue $ws.Range('D16') '6.358'
Set-TextValue $ws.Range('E16') '  -1.33%  '
Set-TextValue $ws.Range('D17') '29.457.81'
Set-TextValue $ws.Range('E17') '  +1.71%  '
Set-TextValue $ws.Range('D18') '0.000008283'
Set-TextValue $ws.Range('E18') '  -0.86%  '
Set-TextValue $ws.Range('D19') '2.099.38'
Set-TextValue $ws.Range('E19') '  +0.78%  '
Set-TextValue $ws.Range('D20') '238.14'
Set-TextValue $ws.Range('E20') '  -2.37%  '
Set-TextValue $ws.Range('D21') '12.72'
Set-TextValue $ws.Range('E21') '  -0.21%  '
Set-TextValue $ws.Range('E22') '  +0.02%  '
Set-TextValue $ws.Range('D23') '7.640'
Set-TextValue $ws.Range('E23') '  +2.06%  '
Set-TextValue $ws.Range('E24') '  -0.01%  '
Set-TextValue $ws.Range('D25') '0.1494'
Set-TextValue $ws.Range('E25') '  +1.48%  '
Set-TextValue $ws.Range('D26') '8.902'
Set-TextValue $ws.Range('E26') '  +1.06%  '
Set-TextValue $ws.Range('D27') '159.68'
Set-TextValue $ws.Range('E27') '  -2.27%  '
Set-TextValue $ws.Range('D28') '18.27'
Set-TextValue $ws.Range('E28') '  +0.17%  '
Set-TextValue $ws.Range('D29') '1.533'
Set-TextValue $ws.Range('E29') '  -1.32%  '
Set-TextValue $ws.Range('D30') '4.245'
Set-TextValue $ws.Range('E30') '  +0.57%  '
Set-TextValue $ws.Range('D31') '4.158'
Set-TextValue $ws.Range('E31') '  -0.09%  '
Set-TextValue $ws.Range('D32') '1.213'
Set-TextValue $ws.Range('E32') '  +3.40%  '
Set-TextValue $ws.Range('D33') '0.05102'
Set-TextValue $ws.Range('E33') '  -0.40%  '
Set-TextValue $ws.Range('D34') '0.7721'
Set-TextValue $ws.Range('E34') '  -0.19%  '
Set-TextValue $ws.Range('D35') '1.886'
Set-TextValue $ws.Range('E35') '  +1.96%  '
Set-TextValue $ws.Range('D36') '1.148'
Set-TextValue $ws.Range('E36') '  +0.32%  '
Set-TextValue $ws.Range('D37') '2.675'
Set-TextValue $ws.Range('E37') '  -0.17%  '
Set-TextValue $ws.Range('D38') '1.333.56'
Set-TextValue $ws.Range('E38') '  +7.34%  '
Set-TextValue $ws.Range('D39') '0.01867'
Set-TextValue $ws.Range('E39') '  +0.97%  '
Set-TextValue $ws.Range('E40') '  +0.65%  '
Set-TextValue $ws.Range('D41') '0.9592'
Set-TextValue $ws.Range('E41') '  +2.08%  '
Set-TextValue $ws.Range('B42') 'FraxShare'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D42') '5.842'
Set-TextValue $ws.Range('E42') '  +1.77%  '
Set-TextValue $ws.Range('B43') 'Quant'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D43') '105.71'
Set-TextValue $ws.Range('E43') '  -2.76%  '
Set-TextValue $ws.Range('D44') '1.000'
Set-TextValue $ws.Range('E44') '  +0.08%  '
Set-TextValue $ws.Range('D45') '9.816'
Set-TextValue $ws.Range('E45') '  +2.13%  '
Set-TextValue $ws.Range('D46') '0.00000000125'
Set-TextValue $ws.Range('E46') '  +2.05%  '
Set-TextValue $ws.Range('D47') '1.998.19'
Set-TextValue $ws.Range('E47') '  +0.74%  '
Set-TextValue $ws.Range('D48') '0.5216'
Set-TextValue $ws.Range('E48') '  +0.85%  '
Set-TextValue $ws.Range('D49') '1.783'
Set-TextValue $ws.Range('E49') '  +1.93%  '
Set-TextValue $ws.Range('D50') '63.13'
Set-TextValue $ws.Range('E50') '  -1.82%  '
Set-TextValue $ws.Range('D51') '6.970'
Set-TextValue $ws.Range('E51') '  +0.56%  '
